$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column (C) for the modules that are now enabled ("Y")
$ws.Range("C2:C6").Value = "Y"

# Move the active selection to G7, matching the saved sheet view
$ws.Range("G7").Select()
